# Natmi following Dr Hou advice
# Update the LR-pair result sheet with recomputed NatMI statistics for
# Ccl11-Ackr2 now that 3 cells (instead of 1) are counted as expressing
# the ligand/receptor, and downstream derived metrics are recalculated
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs/sCs, Ccl11-Ackr2)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7503183333333333
$ws.Range("H2").Value = 2.250955
$ws.Range("I2").Value = 0.001871730106429624
$ws.Range("J2").Value = 0.001871730106429624
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.35607266666667
$ws.Range("N2").Value = 37.068218
$ws.Range("Q2").Value = 9.270987849798889
$ws.Range("R2").Value = 83.43889064819
$ws.Range("S2").Value = 0.001871730106429624
$ws.Range("T2").Value = 0.001871730106429624

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 383.1307676666667
$ws.Range("H3").Value = 1149.392303
$ws.Range("I3").Value = 0.9557508602453543
$ws.Range("J3").Value = 0.9557508602453542
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.35607266666667
$ws.Range("N3").Value = 37.068218
$ws.Range("Q3").Value = 4733.991606125118
$ws.Range("R3").Value = 42605.92445512606
$ws.Range("S3").Value = 0.9557508602453543
$ws.Range("T3").Value = 0.9557508602453542

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.98778433333333
$ws.Range("H4").Value = 50.963353
$ws.Range("I4").Value = 0.04237740964821621
$ws.Range("J4").Value = 0.0423774096482162
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.35607266666667
$ws.Range("N4").Value = 37.068218
$ws.Range("Q4").Value = 209.9022976683283
$ws.Range("R4").Value = 1889.120679014954
$ws.Range("S4").Value = 0.04237740964821621
$ws.Range("T4").Value = 0.0423774096482162
